# SH11SIMPOP_StockAge.xlsx - rename stock-subgroup column headers (row 1,
# columns I:BF on Sheet1) from the old "BY<yy><STOCK>" convention to the
# new "<STOCK>..BY<yy>" convention (double-dot separator), per commit
# "Fix things to work with stock-by-sex and stock-by-age".
#
# Old names look like BY04UPSALM, BY05UPSALM, ... BY08LSNAKE
# New names look like UPSALM..BY04, UPSALM..BY05, ... LSNAKE..BY08

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value  = "UPSALM..BY04"
$ws.Range("J1").Value  = "UPSALM..BY05"
$ws.Range("K1").Value  = "UPSALM..BY06"
$ws.Range("L1").Value  = "UPSALM..BY07"
$ws.Range("M1").Value  = "UPSALM..BY08"
$ws.Range("N1").Value  = "MFSALM..BY04"
$ws.Range("O1").Value  = "MFSALM..BY05"
$ws.Range("P1").Value  = "MFSALM..BY06"
$ws.Range("Q1").Value  = "MFSALM..BY07"
$ws.Range("R1").Value  = "MFSALM..BY08"
$ws.Range("S1").Value  = "SFSALM..BY04"
$ws.Range("T1").Value  = "SFSALM..BY05"
$ws.Range("U1").Value  = "SFSALM..BY06"
$ws.Range("V1").Value  = "SFSALM..BY07"
$ws.Range("W1").Value  = "SFSALM..BY08"
$ws.Range("X1").Value  = "LOSALM..BY04"
$ws.Range("Y1").Value  = "LOSALM..BY05"
$ws.Range("Z1").Value  = "LOSALM..BY06"
$ws.Range("AA1").Value = "LOSALM..BY07"
$ws.Range("AB1").Value = "LOSALM..BY08"
$ws.Range("AC1").Value = "UPCLWR..BY04"
$ws.Range("AD1").Value = "UPCLWR..BY05"
$ws.Range("AE1").Value = "UPCLWR..BY06"
$ws.Range("AF1").Value = "UPCLWR..BY07"
$ws.Range("AG1").Value = "UPCLWR..BY08"
$ws.Range("AH1").Value = "SFCLWR..BY04"
$ws.Range("AI1").Value = "SFCLWR..BY05"
$ws.Range("AJ1").Value = "SFCLWR..BY06"
$ws.Range("AK1").Value = "SFCLWR..BY07"
$ws.Range("AL1").Value = "SFCLWR..BY08"
$ws.Range("AM1").Value = "LOCLWR..BY04"
$ws.Range("AN1").Value = "LOCLWR..BY05"
$ws.Range("AO1").Value = "LOCLWR..BY06"
$ws.Range("AP1").Value = "LOCLWR..BY07"
$ws.Range("AQ1").Value = "LOCLWR..BY08"
$ws.Range("AR1").Value = "IMNAHA..BY04"
$ws.Range("AS1").Value = "IMNAHA..BY05"
$ws.Range("AT1").Value = "IMNAHA..BY06"
$ws.Range("AU1").Value = "IMNAHA..BY07"
$ws.Range("AV1").Value = "IMNAHA..BY08"
$ws.Range("AW1").Value = "GRROND..BY04"
$ws.Range("AX1").Value = "GRROND..BY05"
$ws.Range("AY1").Value = "GRROND..BY06"
$ws.Range("AZ1").Value = "GRROND..BY07"
$ws.Range("BA1").Value = "GRROND..BY08"
$ws.Range("BB1").Value = "LSNAKE..BY04"
$ws.Range("BC1").Value = "LSNAKE..BY05"
$ws.Range("BD1").Value = "LSNAKE..BY06"
$ws.Range("BE1").Value = "LSNAKE..BY07"
$ws.Range("BF1").Value = "LSNAKE..BY08"

# Cosmetic view state matching the re-saved workbook: zoomed to 150%,
# with the last active selection at BI12 (off the used range, as in the
# target sheetView).
$win = $excel.ActiveWindow
$win.Zoom = 150
$ws.Range("BI12").Select()
